# Applies updated "想去人数" (want-to-go count) values across the
# 展览 (sheet1), 演出 (sheet2), and 全部类型 (sheet4) worksheets,
# matching the regenerated gh-pages data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    "F5"  = 4533
    "F8"  = 125
    "F9"  = 3060
    "F11" = 580
    "F13" = 585
    "F15" = 506
    "F16" = 352
    "F17" = 128
    "F18" = 1754
    "F19" = 1292
    "F21" = 1534
    "F24" = 40
    "F27" = 36
    "F28" = 85
    "F29" = 118
    "F31" = 3445
    "F32" = 739
    "F33" = 60
    "F34" = 225
    "F35" = 53
    "F36" = 1671
}
foreach ($addr in $sheet1Updates.Keys) {
    $ws1.Range($addr).Value = $sheet1Updates[$addr]
}

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 34

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    "F5"  = 4533
    "F8"  = 125
    "F9"  = 3060
    "F11" = 580
    "F13" = 585
    "F15" = 506
    "F17" = 352
    "F18" = 128
    "F19" = 1754
    "F20" = 1292
    "F22" = 1534
    "F25" = 40
    "F28" = 36
    "F29" = 85
    "F30" = 118
    "F32" = 3445
    "F33" = 34
    "F34" = 739
    "F35" = 60
    "F36" = 225
    "F37" = 53
    "F38" = 1671
}
foreach ($addr in $sheet4Updates.Keys) {
    $ws4.Range($addr).Value = $sheet4Updates[$addr]
}
